$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns
$ws.Range("I1").Value = "checkOutDate"
$ws.Range("J1").Value = "totalBill"

# Update existing checkInDate value for row 2
$ws.Range("G2").Value = "2021-10-04T18:18:33.152Z"

# Add new data for row 2
$ws.Range("I2").Value = "2021-10-04T19:24:43.416Z"
$ws.Range("J2").Value = 767
